$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.157.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.84%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.913.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.58%  "

$ws.Range("E4").Value = "  -1.42%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.16%  "

$ws.Range("E7").Value = "  -5.82%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4010"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.61%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.22"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08385"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.71%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.044"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.98%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.889.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -8.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.420"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.14%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.058"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.24%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001066"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06606"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.753"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.136.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.306"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.117.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.94%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.34%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.763"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.89%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.129"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.05%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9756"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.79%  "

$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09646"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.443"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.58%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.553"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.41%  "

$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.630"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.273"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.812"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.93%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02296"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06146"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6163"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.29%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1908"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.17%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.300"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5866"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.39%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.75"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.94%  "

$ws.Range("E48").Value = "  -5.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.436"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06912"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.21%  "
